$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the left table by one more year column (2023), matching the
# formatting already used by the 2022 column (J).
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1135.9000000000001
$ws.Range("K5").Value = 970.3
$ws.Range("K6").Value = 1234.9000000000001
